# "Checkliste und Todo Update"
# - Update D18 (row "Stackfunktionen (incl. Visualisierung)") from 0 to 2 points.
#   The SUM formulas in D34/D35 recalc automatically from this change.
# - Update the saved view/selection state of the sheet: scroll position and
#   the active selected cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")
$ws.Activate()

# Update the IST points for "Stackfunktionen (incl. Visualisierung)"
$ws.Range("D18").Value = 2

# Restore the view: scroll so A6 is the top-left visible cell, and select D13
$excel.ActiveWindow.ScrollRow = 6
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D13").Select()
